# transactionInBank/checklist.xlsx -- "bank application with checks and junit"
#
# Row 7 ("4.junit") gets marked "done" for the createBank() / getBankDetails()
# / createCustomer() columns (C, D, E), matching the pattern already used by
# the other checklist rows. The active selection also moves from O11 to A10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = "done"
$ws.Range("D7").Value = "done"
$ws.Range("E7").Value = "done"

# Update the saved cursor/selection to A10 (was O11).
[void]$ws.Range("A10").Select()
